$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Replace the text of a single run (found verbatim, case-sensitive, whole match)
# with new text, without disturbing any sibling runs in the paragraph (avoids
# the engine's run auto-coalescing that a plain Find/Replace triggers).
function Replace-RunText($old, $new) {
    $searchRng = $d.Content
    $found = $searchRng.Find.Execute(
        $old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0
    )
    if (-not $found) {
        throw "Could not find text: $old"
    }
    # Re-seat a clean Range over the same span - calling InsertXML directly on
    # the Find-narrowed range duplicates content instead of replacing it.
    $targetRng = $d.Range($searchRng.Start, $searchRng.End)
    $escaped = $new -replace '&', '&amp;'
    $xml = "<w:p $wNs><w:r><w:t xml:space=`"preserve`">$escaped</w:t></w:r></w:p>"
    $targetRng.InsertXML($xml)
}

# 1. International Payments bullet - add "increasing revenue by 15%."
Replace-RunText `
    "existing Domestic payment flows with International options, to enable users to make International Payments with current exchange rates." `
    "existing Domestic payment flows with International options, to enable users to make International Payments with current exchange rates, increasing revenue by 15%."

# 2. Local Environment Stability bullet - shorten / reword ending
Replace-RunText `
    "Local Environment Stability Issues with a self-directed NodeJS project that automated engineering tasks, including; cloning 11 repos pointing to multiple remotes per brand, running 4-5 server commands, creating multiple feature branches, updating local branches to the latest release code," `
    "Local Environment Stability Issues with a self-directed NodeJS project that automated engineering tasks, including; cloning 11 repos for multiple remotes per brand. This reduced the first-time setup for local development from about 1 week to 1/2 a day."

# 3. Design guide bullet - add "reducing lines of code by 25%."
Replace-RunText `
    "a design guide with classNames and partials, from the product design system, to build stylesheets which implemented reusable classnames and styles." `
    "a design guide with classNames and partials, from the product design system, to build stylesheets which implemented reusable classnames and styles, reducing lines of code by 25%."

# 4. Responsive Design bullet - "lift" becomes "increased ad engagement by 30%"
Replace-RunText `
    "Responsive Design solution that consistently generated lift across all existing playable ads, and became an engineering standard." `
    "Responsive Design solution that consistently generated increased ad engagement by 30% across all existing playable ads, and became an engineering standard."

# 5. Insert a brand-new "Optimized" bullet before the "Managed & Migrated" bullet,
#    and extend the "Managed & Migrated" bullet text.

# Locate the "Managed & Migrated pull requests" bullet paragraph.
$targetIndex = -1
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "Managed*Migrated pull requests*") {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq -1) {
    throw "Could not locate the 'Managed & Migrated pull requests' paragraph"
}

$managedPara = $d.Paragraphs.Item($targetIndex)
$managedRange = $managedPara.Range
$managedRange.InsertParagraphBefore()

# The new blank paragraph is now immediately before the (shifted) Managed paragraph.
$newPara = $d.Paragraphs.Item($targetIndex)
$newRange = $newPara.Range

$newParaXml = "<w:p $wNs>" +
    '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1002"/></w:numPr><w:pStyle w:val="Compact"/></w:pPr>' +
    '<w:r><w:rPr><w:bCs/><w:b/></w:rPr><w:t xml:space="preserve">Optimized</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">feature branch creation for our 4 repos with my script that created multiple feature branches with the JIRA ticket number, and updated local branches to the latest release code.</w:t></w:r>' +
    '</w:p>'
$newRange.InsertXML($newParaXml)

# Now extend the Managed paragraph's trailing sentence (run contains a literal "&").
Replace-RunText `
    "& Migrated pull requests from team-specific repos, to Bedrock & Release Management repos, including batch cherry-picked commits from my team to the upstream repos." `
    "& Migrated pull requests from team-specific repos, to Bedrock & Release Management repos, including batch cherry-picked commits from my team to the upstream repos. I was responsible for making sure all of the features for the release were included in the release branches."
